# Apply hybrid bold + color (2C3E50) highlighting to quantitative impact
# metrics (percentages, dollar amounts, large numbers) across the resume's
# experience bullets and achievements, without altering the underlying text.

function Set-MetricHighlight($paraRange, $text) {
    $r = $paraRange.Duplicate
    $ok = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($ok) {
        $r.Font.Bold = 1
        $r.Font.Color = 5258796   # RGB(0x2C,0x3E,0x50) == w:color 2C3E50
    }
    return $ok
}

$d = $word.ActiveDocument

# 1) Siege Analytics - race coding discovery bullet: "... from 23% to 64%"
$p = $d.Paragraphs.Item(9).Range
Set-MetricHighlight $p '23%' | Out-Null
Set-MetricHighlight $p '64%' | Out-Null

# 2) Siege Analytics - turnout prediction bullet: "87% ... 71% ... ±4.2% to ±2.1%"
$p = $d.Paragraphs.Item(11).Range
Set-MetricHighlight $p '87%' | Out-Null
Set-MetricHighlight $p '71%' | Out-Null
Set-MetricHighlight $p '±4.2%' | Out-Null
Set-MetricHighlight $p '±2.1%' | Out-Null

# 3) Myers Research - RFP bullet: "bids from 1,200 vendors"
$p = $d.Paragraphs.Item(31).Range
Set-MetricHighlight $p '1,200' | Out-Null

# 4) Lake Research Partners - Polling Consortium bullet: "$400M ... $1B+"
$p = $d.Paragraphs.Item(46).Range
Set-MetricHighlight $p '$400M' | Out-Null
Set-MetricHighlight $p '$1B' | Out-Null

# 5) Key Achievements - Algorithm bullet: "73.5% ... $4.7M"
$p = $d.Paragraphs.Item(63).Range
Set-MetricHighlight $p '73.5%' | Out-Null
Set-MetricHighlight $p '$4.7M' | Out-Null

# 6) Key Achievements - turnout prediction bullet (no trailing margin clause): "87% ... 71%"
$p = $d.Paragraphs.Item(65).Range
Set-MetricHighlight $p '87%' | Out-Null
Set-MetricHighlight $p '71%' | Out-Null
